try {
    $wb = $excel.ActiveWorkbook
    $ws = $wb.ActiveSheet

    # The first block of task rows (A2:I6 -- Task IDs 0..4) represents one
    # "routing" cycle. The edit combines routing and scheduling by repeating
    # that block twice more, continuing the Task ID sequence (10..19) while
    # keeping every other column identical to the source rows.
    $ws.Range("A2:I6").Copy()
    $ws.Range("A12:I16").Insert()

    $ws.Range("A2:I6").Copy()
    $ws.Range("A17:I21").Insert()

    # Fix up the Task ID column (A) for the newly inserted rows so it
    # continues the existing sequence instead of repeating 0..4 twice.
    for ($i = 0; $i -lt 10; $i++) {
        $ws.Cells.Item(12 + $i, 1).Value = 10 + $i
    }

    # Move the active selection to match the saved view state.
    $ws.Range("B24").Select()
} catch {
    Write-Output "ERROR: $_"
}
